$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use Text format on D:E columns while writing so that numeric-looking /
# percent-looking strings are kept as literal text (matching the source
# workbook's inlineStr cell type) instead of being auto-converted by Excel
# into numbers / percentages.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "308.84"
$ws.Range("E2").Value = "0.32%"
$ws.Range("D3").Value = "40.92"
$ws.Range("E3").Value = "0.02%"
$ws.Range("D4").Value = "5.125"
$ws.Range("E4").Value = "1.77%"
$ws.Range("D5").Value = "0.07618"
$ws.Range("E5").Value = "-0.24%"
$ws.Range("D6").Value = "1.604"
$ws.Range("E6").Value = "-0.17%"
$ws.Range("D8").Value = "0.9081"
$ws.Range("E8").Value = "-0.16%"
$ws.Range("D9").Value = "0.1270"
$ws.Range("E9").Value = "25.06%"
$ws.Range("D10").Value = "0.1807"
$ws.Range("E10").Value = "2.31%"
$ws.Range("D11").Value = "0.09097"
$ws.Range("E11").Value = "0.66%"
$ws.Range("D12").Value = "0.04293"
$ws.Range("E12").Value = "-2.47%"
$ws.Range("E13").Value = "-0.81%"
$ws.Range("D14").Value = "0.001254"
$ws.Range("E14").Value = "0.16%"
$ws.Range("D15").Value = "0.005815"
$ws.Range("E15").Value = "0.35%"
$ws.Range("D16").Value = "3.352"
$ws.Range("E16").Value = "-0.40%"
$ws.Range("D17").Value = "4.283"
$ws.Range("E17").Value = "1.16%"
$ws.Range("D19").Value = "6.929"
$ws.Range("E19").Value = "2.25%"
$ws.Range("D20").Value = "0.1394"
$ws.Range("E20").Value = "2.80%"
$ws.Range("E21").Value = "-1.25%"
$ws.Range("D22").Value = "0.04036"
$ws.Range("E22").Value = "-2.73%"
$ws.Range("E23").Value = "4.88%"
$ws.Range("D24").Value = "0.004042"
$ws.Range("E24").Value = "-1.23%"
$ws.Range("E25").Value = "-2.10%"
$ws.Range("E26").Value = "24.81%"
$ws.Range("D38").Value = "0.02416"
$ws.Range("E38").Value = "0.11%"
$ws.Range("D39").Value = "0.05213"
$ws.Range("E39").Value = "0.87%"
$ws.Range("D40").Value = "0.007839"
$ws.Range("E40").Value = "0.95%"
$ws.Range("D41").Value = "0.1301"
$ws.Range("E41").Value = "-0.66%"
$ws.Range("D42").Value = "0.006807"
$ws.Range("E42").Value = "-3.83%"
$ws.Range("D43").Value = "0.001953"
$ws.Range("E43").Value = "0.32%"
$ws.Range("D44").Value = "0.007350"
$ws.Range("E44").Value = "-8.47%"
$ws.Range("D45").Value = "0.3371"
$ws.Range("E45").Value = "10.23%"
$ws.Range("D46").Value = "0.00006890"
$ws.Range("E46").Value = "8.00%"
$ws.Range("E47").Value = "0.23%"
$ws.Range("D48").Value = "0.1399"
$ws.Range("E48").Value = "2,184.95%"
$ws.Range("E49").Value = "-31.72%"
$ws.Range("E50").Value = "0.23%"
$ws.Range("E51").Value = "0.23%"

# Restore the default ("Normal") style on the touched range so the cells do
# not carry a leftover text-format style index (matches original workbook
# where these cells had no explicit style).
$ws.Range("D2:E51").Style = "Normal"
